$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's format onto the two new header cells,
# then set their text (mirrors the existing header cells H1.. style, s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data column values (rows 2-5), unstyled like the other data cells
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 7
